# "Working on t test"
# Clear out the stub placeholder values (1/1/1/1/1/1, 2/2/2/2/2/2 and the
# "0 (P=0.050)" placeholder text) that used to sit in the header/summary
# row of each table block. Rows 9 and 48 additionally get a real computed
# statistic (0.6133371569099184) written into the "Treatment at T2" (I:N)
# columns, while their "Control at T1" (C:H) and P-value (O:T) columns are
# cleared just like the rest.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose entire C:T block just gets wiped back to blank.
$fullClearRows = @(3, 13, 19, 25, 31, 37, 42, 52, 58, 64, 70, 76)

foreach ($r in $fullClearRows) {
    $ws.Range("C$r`:T$r").ClearContents()
}

# Rows that keep a real value in I:N, but still clear C:H and O:T.
$partialRows = @(9, 48)

foreach ($r in $partialRows) {
    $ws.Range("C$r`:H$r").ClearContents()
    $ws.Range("I$r`:N$r").Value = 0.6133371569099184
    $ws.Range("O$r`:T$r").ClearContents()
}
